# Weekly update: insert a new week's worth of data (5 rows) for
# "Fruta, Terminal La Palmera de La Serena - Durazno" ahead of the
# existing history, pushing the previous rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 361; this shifts the existing
# rows 361-376 down to 366-381 and grows the used range accordingly.
$ws.Rows.Item(361).Resize(5).Insert()

# Columns that are constant across every data row in this subset.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria   = "Durazno"
$unidadBins  = "`$/bins (400 kilos)"

function Set-Row {
    param($r, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg, $kgUnidad)

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidadBins
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 361 44610 "Doctor Davis"   "Especial" 16 395000 400000 397500 "Región de O'Higgins" 994 400
Set-Row 362 44610 "Doctor Davis"   "Primera"  16 365000 370000 367500 "Región de O'Higgins" 919 400
Set-Row 363 44610 "Doctor Davis"   "Segunda"  10 325000 330000 327500 "Región de O'Higgins" 819 400
Set-Row 364 44610 "September Sun"  "Especial" 10 385000 390000 387500 "Región de O'Higgins" 969 400
Set-Row 365 44610 "September Sun"  "Primera"  10 335000 340000 337500 "Región de O'Higgins" 844 400

# Give the new date cells the same date/time number format as the rest
# of column D.
$ws.Range("D361:D365").NumberFormat = $ws.Range("D366").NumberFormat
